$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.313.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.577.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +2.66%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("E6").Value = "  -0.56%  "

$ws.Range("E7").Value = "  +2.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.60%  "

$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0881"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.802.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.575.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.341.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.04%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0699"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("E25").Value = "  +4.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("E32").Value = "  -1.06%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.385.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.93%  "

$ws.Range("E38").Value = "  +3.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.26%  "

$ws.Range("E40").Value = "  -1.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.537"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.77%  "

$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("E43").Value = "  +3.01%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.24%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.982"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.713.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
